$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.521.88'
$ws.Range('E2').Value = '  +0.03%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.558.13'
$ws.Range('E3').Value = '  +3.55%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.39'
$ws.Range('E5').Value = '  +1.94%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.87'
$ws.Range('E6').Value = '  +3.75%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.559.15'
$ws.Range('E7').Value = '  +3.59%  '

$ws.Range('E8').Value = '  +0.07%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('E9').Value = '  +0.21%  '

$ws.Range('E10').Value = '  +2.98%  '

$ws.Range('E11').Value = '  -4.15%  '

$ws.Range('E12').Value = '  +4.76%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.159.68'
$ws.Range('E13').Value = '  +3.56%  '

$ws.Range('E14').Value = '  +4.03%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.26'
$ws.Range('E15').Value = '  +2.45%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.556.31'
$ws.Range('E16').Value = '  +2.98%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.117'
$ws.Range('E17').Value = '  +1.65%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '65.529.79'
$ws.Range('E18').Value = '  +0.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.38'
$ws.Range('E19').Value = '  +4.99%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.93'
$ws.Range('E20').Value = '  +1.95%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.33'
$ws.Range('E21').Value = '  +5.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '395.74'
$ws.Range('E22').Value = '  +0.24%  '

$ws.Range('E23').Value = '  +4.87%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.697.24'
$ws.Range('E24').Value = '  +3.31%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.21'
$ws.Range('E25').Value = '  +1.23%  '

$ws.Range('E26').Value = '  +0.08%  '

$ws.Range('E27').Value = '  +11.05%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.88'
$ws.Range('E28').Value = '  +8.96%  '

$ws.Range('E29').Value = '  +0.03%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.31'
$ws.Range('E30').Value = '  +1.98%  '

$ws.Range('E31').Value = '  +3.00%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.570.26'
$ws.Range('E32').Value = '  +3.73%  '

$ws.Range('E33').Value = '  +0.06%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.148'
$ws.Range('E34').Value = '  +0.16%  '

$ws.Range('E35').Value = '  +3.72%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.29'
$ws.Range('E36').Value = '  +6.31%  '

$ws.Range('E37').Value = '  +2.08%  '

$ws.Range('E38').Value = '  +3.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '167.56'
$ws.Range('E39').Value = '  -2.30%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.04'
$ws.Range('E40').Value = '  +4.78%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0812'
$ws.Range('E41').Value = '  +5.37%  '

$ws.Range('E42').Value = '  +1.45%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.41'
$ws.Range('E43').Value = '  +16.43%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '43.13'
$ws.Range('E44').Value = '  -0.83%  '

$ws.Range('E45').Value = '  -0.09%  '

$ws.Range('E46').Value = '  +0.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.71'
$ws.Range('E47').Value = '  +5.97%  '

$ws.Range('E48').Value = '  +9.62%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.488.90'
$ws.Range('E49').Value = '  +13.07%  '

$ws.Range('E50').Value = '  +4.15%  '

$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.35'
$ws.Range('E51').Value = '  +18.80%  '
